# Update gh-pages to output generated at 456a3b4
# Updates the "want to go" count column (F) on several sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 625
$wsExhibition.Range("F8").Value = 1255
$wsExhibition.Range("F9").Value = 3941
$wsExhibition.Range("F10").Value = 84

# Sheet "演出" (Performances)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 53

# Sheet "全部类型" (All Types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 625
$wsAll.Range("F8").Value = 1256
$wsAll.Range("F9").Value = 3941
$wsAll.Range("F10").Value = 84
$wsAll.Range("F11").Value = 53
